$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header in H1, copying the header formatting (bold,
# centered, bordered) from the neighboring "sum" header in G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for each data row.
$ws.Range("H2:H4").Value = 0
$ws.Range("H5").Value = 1
